$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2848.4285
$ws.Range("I6").Value = 2848.4285
$ws.Range("K6").Value = 8545.2855
$ws.Range("M6").Value = -8433.2855
$ws.Range("H8").Value = 151
$ws.Range("I8").Value = 202
$ws.Range("K8").Value = 606
$ws.Range("M8").Value = -467
$ws.Range("H17").Value = 607836.0600000001
$ws.Range("I17").Value = 1080.3
$ws.Range("J17").Value = 871642.9399999999
$ws.Range("K17").Value = 3240.9
$ws.Range("L17").Value = 2614928.82
$ws.Range("M17").Value = -3072.9
$ws.Range("N17").Value = -2615264.82
$ws.Range("H41").Value = 2201.125
$ws.Range("I41").Value = 2351.5
$ws.Range("K41").Value = 2351.5
$ws.Range("M41").Value = -1911.5
$ws.Range("H58").Value = 2534
$ws.Range("I58").Value = 3126
$ws.Range("J58").Value = 1350
$ws.Range("K58").Value = 9378
$ws.Range("L58").Value = 4050
$ws.Range("M58").Value = -9228
$ws.Range("N58").Value = -4350
$ws.Range("H64").Value = 7363.75
$ws.Range("I64").Value = 3303
$ws.Range("J64").Value = 9800.200000000001
$ws.Range("K64").Value = 3303
$ws.Range("L64").Value = 9800.200000000001
$ws.Range("M64").Value = -3055
$ws.Range("N64").Value = -10296.2
$ws.Range("H67").Value = 7363.75
$ws.Range("I67").Value = 3303
$ws.Range("J67").Value = 9800.200000000001
$ws.Range("K67").Value = 3303
$ws.Range("L67").Value = 9800.200000000001
$ws.Range("M67").Value = -2445
$ws.Range("N67").Value = -11516.2
$ws.Range("H76").Value = 2994.5557
$ws.Range("I76").Value = 2994.5557
$ws.Range("K76").Value = 2994.5557
$ws.Range("M76").Value = -2679.5557
$ws.Range("H79").Value = 2994.5557
$ws.Range("I79").Value = 2994.5557
$ws.Range("K79").Value = 2994.5557
$ws.Range("M79").Value = -1902.5557
$ws.Range("H121").Value = 4487.385
$ws.Range("J121").Value = 4487.385
$ws.Range("L121").Value = 13462.155
$ws.Range("N121").Value = -16956.155
$ws.Range("H141").Value = 5859.9
$ws.Range("I141").Value = 5844.8887
$ws.Range("J141").Value = 5995
$ws.Range("K141").Value = 17534.6661
$ws.Range("L141").Value = 17985
$ws.Range("M141").Value = -12354.6661
$ws.Range("N141").Value = -28345

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 2515002.5
$ws.Range("I8").Value = 2515002.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 2515002.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -2514858.5
$ws.Range("N8").ClearContents()
$ws.Range("H32").Value = 2340.3137
$ws.Range("I32").Value = 2007.1459
$ws.Range("K32").Value = 2007.1459
$ws.Range("M32").Value = -1720.1459
$ws.Range("H45").Value = 4969
$ws.Range("I45").Value = 4276.7
$ws.Range("K45").Value = 4276.7
$ws.Range("M45").Value = -3899.7
$ws.Range("H74").Value = 2153.8333
$ws.Range("I74").Value = 2077
$ws.Range("K74").Value = 2077
$ws.Range("M74").Value = -1203
$ws.Range("H77").Value = 2153.8333
$ws.Range("I77").Value = 2077
$ws.Range("K77").Value = 10385
$ws.Range("M77").Value = -6017
$ws.Range("H132").Value = 2533.4375
$ws.Range("I132").Value = 2291.9534
$ws.Range("J132").Value = 4610.2
$ws.Range("K132").Value = 6875.860199999999
$ws.Range("L132").Value = 13830.6
$ws.Range("M132").Value = -4345.860199999999
$ws.Range("N132").Value = -18890.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3149.3
$ws.Range("I20").Value = 2414.8333
$ws.Range("J20").Value = 4251
$ws.Range("K20").Value = 2414.8333
$ws.Range("L20").Value = 4251
$ws.Range("M20").Value = -2167.8333
$ws.Range("N20").Value = -4745
$ws.Range("H56").Value = 23533.334
$ws.Range("I56").Value = 17000
$ws.Range("J56").Value = 24000
$ws.Range("K56").Value = 17000
$ws.Range("L56").Value = 24000
$ws.Range("M56").Value = -16261
$ws.Range("N56").Value = -25478

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3358.1333
$ws.Range("J58").Value = 3648.4736
$ws.Range("L58").Value = 3648.4736
$ws.Range("N58").Value = -4054.4736
$ws.Range("H132").Value = 2745
$ws.Range("I132").Value = 2314.5
$ws.Range("J132").Value = 3606
$ws.Range("K132").Value = 6943.5
$ws.Range("L132").Value = 10818
$ws.Range("M132").Value = -4413.5
$ws.Range("N132").Value = -15878
$ws.Range("H136").Value = 3358.1333
$ws.Range("J136").Value = 3648.4736
$ws.Range("L136").Value = 10945.4208
$ws.Range("N136").Value = -16045.4208

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 117.454544
$ws.Range("J2").Value = 134.28572
$ws.Range("L2").Value = 805.71432
$ws.Range("N2").Value = -1031.71432
$ws.Range("H9").Value = 1800
$ws.Range("I9").Value = 1800
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 5400
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -5176
$ws.Range("N9").ClearContents()
$ws.Range("H12").Value = 35.666668
$ws.Range("J12").Value = 48.8
$ws.Range("L12").Value = 146.4
$ws.Range("N12").Value = -492.4
$ws.Range("H61").Value = 472.2
$ws.Range("J61").Value = 539.25
$ws.Range("L61").Value = 1617.75
$ws.Range("N61").Value = -2047.75
$ws.Range("H131").Value = 1635.3036
$ws.Range("I131").Value = 1141.6666
$ws.Range("J131").Value = 1815.9025
$ws.Range("K131").Value = 3424.9998
$ws.Range("L131").Value = 5447.7075
$ws.Range("M131").Value = 1615.0002
$ws.Range("N131").Value = -15527.7075

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1696.8334
$ws.Range("I113").Value = 1798
$ws.Range("K113").Value = 1798
$ws.Range("M113").Value = 372
$ws.Range("H126").Value = 3762.889
$ws.Range("I126").Value = 2998.2222
$ws.Range("K126").Value = 8994.6666
$ws.Range("M126").Value = -6524.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 1500
$ws.Range("I12").Value = 1500
$ws.Range("K12").Value = 1500
$ws.Range("M12").Value = -1330
$ws.Range("H40").Value = 5060.2354
$ws.Range("I40").Value = 4358.9287
$ws.Range("K40").Value = 4358.9287
$ws.Range("M40").Value = -4222.9287
$ws.Range("H46").Value = 2794.2942
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 3073.5334
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 3073.5334
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -3449.5334
$ws.Range("H68").Value = 3235.7273
$ws.Range("I68").Value = 3136.625
$ws.Range("K68").Value = 3136.625
$ws.Range("M68").Value = -2387.625
$ws.Range("H71").Value = 3235.7273
$ws.Range("I71").Value = 3136.625
$ws.Range("K71").Value = 15683.125
$ws.Range("M71").Value = -11939.125
$ws.Range("H82").Value = 2061.0625
$ws.Range("I82").Value = 1941
$ws.Range("K82").Value = 1941
$ws.Range("M82").Value = -1580
$ws.Range("H85").Value = 2061.0625
$ws.Range("I85").Value = 1941
$ws.Range("K85").Value = 1941
$ws.Range("M85").Value = -693
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H122").Value = 30082.455
$ws.Range("I122").Value = 30546.934
$ws.Range("J122").Value = 29087.143
$ws.Range("K122").Value = 91640.802
$ws.Range("L122").Value = 87261.429
$ws.Range("M122").Value = -89190.802
$ws.Range("N122").Value = -92161.429
$ws.Range("H132").Value = 3056.1667
$ws.Range("I132").Value = 2601.3333
$ws.Range("K132").Value = 7803.999899999999
$ws.Range("M132").Value = -5273.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1708.2727
$ws.Range("I136").Value = 1205.561
$ws.Range("K136").Value = 3616.683
$ws.Range("M136").Value = -1066.683
